$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 40000332
$ws.Range("I33").Value = 43478604
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 43478604
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = -43478375
$ws.Range("N33").Value = -658
$ws.Range("H40").Value = 1792.3077
$ws.Range("I40").Value = 1600
$ws.Range("J40").Value = 1912.5
$ws.Range("K40").Value = 1600
$ws.Range("L40").Value = 1912.5
$ws.Range("M40").Value = -1425
$ws.Range("N40").Value = -2262.5
$ws.Range("H64").Value = 3470
$ws.Range("I64").Value = 3450
$ws.Range("J64").Value = 3487.5
$ws.Range("K64").Value = 3450
$ws.Range("L64").Value = 3487.5
$ws.Range("M64").Value = -3202
$ws.Range("N64").Value = -3983.5
$ws.Range("H67").Value = 3470
$ws.Range("I67").Value = 3450
$ws.Range("J67").Value = 3487.5
$ws.Range("K67").Value = 3450
$ws.Range("L67").Value = 3487.5
$ws.Range("M67").Value = -2592
$ws.Range("N67").Value = -5203.5
$ws.Range("H96").Value = 981.5833
$ws.Range("J96").Value = 918.6667
$ws.Range("L96").Value = 2756.0001
$ws.Range("N96").Value = -5502.0001
$ws.Range("H124").Value = 60000
$ws.Range("J124").Value = 60000
$ws.Range("L124").Value = 60000
$ws.Range("N124").Value = -69820
$ws.Range("H132").Value = 2664.9
$ws.Range("I132").Value = 2445.5107
$ws.Range("K132").Value = 7336.532099999999
$ws.Range("M132").Value = -4806.532099999999
$ws.Range("H137").Value = 4000852.5
$ws.Range("I137").Value = 814.2273
$ws.Range("K137").Value = 2442.6819
$ws.Range("M137").Value = 107.3181
$ws.Range("H139").Value = 50172.25
$ws.Range("J139").Value = 50172.25
$ws.Range("L139").Value = 50172.25
$ws.Range("N139").Value = -60452.25
$ws.Range("H140").Value = 78940
$ws.Range("J140").Value = 78940
$ws.Range("L140").Value = 78940
$ws.Range("N140").Value = -89300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 29413848
$ws.Range("I61").Value = 38463852
$ws.Range("K61").Value = 38463852
$ws.Range("M61").Value = -38463640
$ws.Range("H74").Value = 13891939
$ws.Range("I74").Value = 22729020
$ws.Range("J74").Value = 5096.7144
$ws.Range("K74").Value = 22729020
$ws.Range("L74").Value = 5096.7144
$ws.Range("M74").Value = -22728146
$ws.Range("N74").Value = -6844.7144
$ws.Range("H77").Value = 13891939
$ws.Range("I77").Value = 22729020
$ws.Range("J77").Value = 5096.7144
$ws.Range("K77").Value = 113645100
$ws.Range("L77").Value = 25483.572
$ws.Range("M77").Value = -113640732
$ws.Range("N77").Value = -34219.572
$ws.Range("H136").Value = 29413848
$ws.Range("I136").Value = 38463852
$ws.Range("K136").Value = 115391556
$ws.Range("M136").Value = -115389006

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3602.3872
$ws.Range("I134").Value = 2212
$ws.Range("J134").Value = 7599.75
$ws.Range("K134").Value = 6636
$ws.Range("L134").Value = 22799.25
$ws.Range("M134").Value = -4101
$ws.Range("N134").Value = -27869.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 20835860
$ws.Range("I132").Value = 25001830
$ws.Range("J132").Value = 6006.5
$ws.Range("K132").Value = 75005490
$ws.Range("L132").Value = 18019.5
$ws.Range("M132").Value = -75002960
$ws.Range("N132").Value = -23079.5
$ws.Range("H141").Value = 114264.586
$ws.Range("J141").Value = 114264.586
$ws.Range("L141").Value = 114264.586
$ws.Range("N141").Value = -124624.586

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 663.75
$ws.Range("I60").Value = 177.5
$ws.Range("J60").Value = 1150
$ws.Range("K60").Value = 532.5
$ws.Range("L60").Value = 3450
$ws.Range("M60").Value = -281.5
$ws.Range("N60").Value = -3952
$ws.Range("H75").Value = 1065.5
$ws.Range("J75").Value = 1034
$ws.Range("L75").Value = 3102
$ws.Range("N75").Value = -5098
$ws.Range("H78").Value = 1065.5
$ws.Range("J78").Value = 1034
$ws.Range("L78").Value = 9306
$ws.Range("N78").Value = -19290
$ws.Range("H121").Value = 1301.5
$ws.Range("J121").Value = 1418.3334
$ws.Range("L121").Value = 4255.0002
$ws.Range("N121").Value = -6875.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4144.625
$ws.Range("I126").Value = 2860.0715
$ws.Range("J126").Value = 5143.722
$ws.Range("K126").Value = 8580.2145
$ws.Range("L126").Value = 15431.166
$ws.Range("M126").Value = -6110.2145
$ws.Range("N126").Value = -20371.166
$ws.Range("H132").Value = 4587.3057
$ws.Range("I132").Value = 3543.75
$ws.Range("K132").Value = 10631.25
$ws.Range("M132").Value = -8101.25
$ws.Range("H138").Value = 59232.832
$ws.Range("J138").Value = 59232.832
$ws.Range("L138").Value = 59232.832
$ws.Range("N138").Value = -69512.83199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1117.8948
$ws.Range("I61").Value = 995.38464
$ws.Range("K61").Value = 995.38464
$ws.Range("M61").Value = -793.38464
$ws.Range("H113").Value = 1117.8948
$ws.Range("I113").Value = 995.38464
$ws.Range("K113").Value = 995.38464
$ws.Range("M113").Value = 1174.61536
$ws.Range("H132").Value = 13167187
$ws.Range("I132").Value = 6777.579
$ws.Range("K132").Value = 20332.737
$ws.Range("M132").Value = -17802.737
$ws.Range("H136").Value = 12504748
$ws.Range("I136").Value = 15626591
$ws.Range("J136").Value = 17375.625
$ws.Range("K136").Value = 46879773
$ws.Range("L136").Value = 52126.875
$ws.Range("M136").Value = -46877223
$ws.Range("N136").Value = -57226.875
$ws.Range("H139").Value = 57543
$ws.Range("J139").Value = 57543
$ws.Range("L139").Value = 57543
$ws.Range("N139").Value = -67823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 36765
$ws.Range("J75").Value = 36765
$ws.Range("L75").Value = 36765
$ws.Range("N75").Value = -38637
$ws.Range("H78").Value = 36765
$ws.Range("J78").Value = 36765
$ws.Range("L78").Value = 110295
$ws.Range("N78").Value = -119655
$ws.Range("H113").Value = 3328.7273
$ws.Range("I113").Value = 2818.25
$ws.Range("J113").Value = 4690
$ws.Range("K113").Value = 8454.75
$ws.Range("L113").Value = 14070
$ws.Range("M113").Value = -6284.75
$ws.Range("N113").Value = -18410
